$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$c = $cs.Item(3)
$c.RGB = 123456
